$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A99").Value = 17
$ws.Range("B99").Value = 15
$ws.Range("C99").Value = 1.5
$ws.Range("D99").Value = 50
$ws.Range("E99").Value = 82.95999999999999
$ws.Range("F99").Value = 10201
